$d = $word.ActiveDocument

# Locate the sentence "duration at different location(by gender)" in
# paragraph "3." using Find, then re-seat a plain Range over the same
# span so InsertXML has a concrete character range to replace.
$finder = $d.Content
$target = "duration at different location(by gender)"
$found = $finder.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find target text '$target'"
}

$rng = $d.Range($finder.Start, $finder.End)

# Replace that span with the expanded wording, reproducing the same
# run/proofErr/bookmark layout Word itself would emit: the text becomes
# "duration at different location(by gender, by age, etc)" and the
# existing "_GoBack" bookmark is carried over between "etc" and ")".
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">duration at different </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>location</w:t></w:r><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>by gender</w:t></w:r><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, by age</w:t></w:r><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>etc</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00E24349"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

Write-Output $d.Content.Text
